$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Iraq League")

# The two match rows (Excel row numbers) whose data got swapped in the
# source spreadsheet (they were re-ordered/re-scraped upstream).
$rowA = 89
$rowB = 90

# All columns whose values belong to the match record itself: B (match id)
# and F..AC (HomeTeam, AwayTeam, score, odds, ...). Columns A (sequence
# number), C, D (league name) and E (date) stay untouched because both
# rows keep their original sequence number / league / date.
$colLetters = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($colLetter in $colLetters) {
    $cellA = $ws.Range($colLetter + $rowA)
    $cellB = $ws.Range($colLetter + $rowB)

    $valueA = $cellA.Value2
    $valueB = $cellB.Value2

    $cellA.Value2 = $valueB
    $cellB.Value2 = $valueA
}
